# Hoàn thiện tool liên quan đến các hình còn lại Circle, Polygon, Oriented Box
#
# This script updates the "name_variable_in_gui" documentation sheet:
#  - Renumbers the STT (sequence) column for rows 18-24 (they had been off by
#    one, duplicating "11").
#  - Removes the now-redundant J/K "label show main/subwindow" helper cells on
#    rows 25-26 (the same screen_main/screen_subwindow_1 pairing is already
#    documented via rows 30-31, columns D:F).
#  - Adds two new GUI-widget rows (24: IMAGE label / label_camera, and
#    25: Combobox3 / btn_resize) used by the newly finished Circle / Polygon /
#    Oriented Box shape tools.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Renumber STT column for rows 18-24 (11..17 -> 12..18) ---
$ws.Range("C18").Value = 12
$ws.Range("C19").Value = 13
$ws.Range("C20").Value = 14
$ws.Range("C21").Value = 15
$ws.Range("C22").Value = 16
$ws.Range("C23").Value = 17
$ws.Range("C24").Value = 18

# --- Row 25: clear STT (no longer used) and the J/K helper note ---
$ws.Range("C25").ClearContents()
$ws.Range("J25").ClearContents()
$ws.Range("K25").ClearContents()

# --- Row 26: clear the J/K helper note (row content in C:F is unchanged) ---
$ws.Range("J26").ClearContents()
$ws.Range("K26").ClearContents()

# --- Apply the same table formatting (font/border) used by the rest of the
#     table to the new spacer row (32) and the two new data rows (33-34) ---
$ws.Range("C24:F24").Copy()
$ws.Range("C32:F34").PasteSpecial(-4122)

# Row 32 stays a blank spacer row.
$ws.Range("C32:F32").ClearContents()

# --- New row 33: IMAGE label ---
# (shared-string table order matters for a byte-faithful round trip: the
# ObjectName "label_camera" was registered before the "IMAGE: " caption, so
# write column E before column D.)
$ws.Range("C33").Value = 24
$ws.Range("E33").Value = "label_camera"
$ws.Range("D33").Value = "IMAGE: "

# --- New row 34: Combobox3 ---
$ws.Range("C34").Value = 25
$ws.Range("D34").Value = "Combobox3"
$ws.Range("E34").Value = "btn_resize"

# --- Update the view's active selection to match the author's cursor ---
$ws.Range("C35").Select()
